$d = $word.ActiveDocument

# Typo / wording correction in the "Ceri has written the descriptive Risk
# assessment ..." bullet: the parenthetical aside about asking Savas for a
# single document is dropped, and the sentence now simply ends with a
# period.
#   "Ceri has written the descriptive Risk assessment (ask Savas if
#    needed in a single document)"
#   becomes
#   "Ceri has written the descriptive Risk assessment."
$rng = $d.Content
$rng.Find.Execute(
    " (ask Savas if needed in a single document)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ".",
    2
)
if (-not $rng.Find.Found) {
    Write-Output "WARNING: target phrase not found"
}
